$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.804.85'
$ws.Range("E2").Value = '  -0.66%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.631.03'
$ws.Range("E3").Value = '  -0.70%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.60'
$ws.Range("E5").Value = '  +0.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5066'
$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2575'
$ws.Range("E8").Value = '  -0.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06422'
$ws.Range("E9").Value = '  +0.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.47'
$ws.Range("E10").Value = '  -2.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07784'
$ws.Range("E11").Value = '  +0.59%  '

$ws.Range("E12").Value = '  -0.62%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.632.44'
$ws.Range("E13").Value = '  -0.65%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.855.53'
$ws.Range("E14").Value = '  -0.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5597'
$ws.Range("E15").Value = '  +2.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅7577'
$ws.Range("E16").Value = '  -2.36%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.95'
$ws.Range("E17").Value = '  -2.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.817.46'
$ws.Range("E18").Value = '  -0.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  +0.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.00'
$ws.Range("E20").Value = '  -1.25%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.325'
$ws.Range("E21").Value = '  -3.30%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.844'
$ws.Range("E22").Value = '  -1.36%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.033'
$ws.Range("E23").Value = '  -1.96%  '

$ws.Range("E24").Value = '  +0.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.800'
$ws.Range("E25").Value = '  -4.89%  '

$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1277'
$ws.Range("E26").Value = '  +1.14%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '140.56'
$ws.Range("E27").Value = '  -1.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.745'
$ws.Range("E28").Value = '  -1.91%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.40'
$ws.Range("E29").Value = '  -1.48%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.237'
$ws.Range("E30").Value = '  -0.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04872'
$ws.Range("E31").Value = '  -0.83%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.284'
$ws.Range("E32").Value = '  +0.43%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.213'
$ws.Range("E33").Value = '  +0.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.549'
$ws.Range("E34").Value = '  -0.18%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.379'
$ws.Range("E35").Value = '  +0.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8970'
$ws.Range("E36").Value = '  -2.43%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.567'
$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.128.20'
$ws.Range("E38").Value = '  -0.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5494'
$ws.Range("E39").Value = '  -1.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01559'
$ws.Range("E40").Value = '  -0.66%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9935'
$ws.Range("E41").Value = '  -0.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.535'
$ws.Range("E42").Value = '  -1.29%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7988'
$ws.Range("E43").Value = '  -0.63%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.21'
$ws.Range("E44").Value = '  -1.59%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.782.63'
$ws.Range("E45").Value = '  +0.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈113'
$ws.Range("E46").Value = '  -5.78%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4440'
$ws.Range("E47").Value = '  -1.95%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.30'
$ws.Range("E48").Value = '  -0.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05051'
$ws.Range("E49").Value = '  -2.61%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.642'
$ws.Range("E50").Value = '  +0.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9997'
$ws.Range("E51").Value = '  -0.17%  '
